# Updated cryptos list with GitHub Actions
# Applies the refreshed prices / 1h-volume percentages scraped for the
# cryptos.xlsx report, plus the corrected Monero / MXToken row ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain text even when it looks like
# a number (e.g. "0.200", "57.50") so Excel doesn't silently coerce it to
# a numeric cell and drop meaningful trailing zeros. We flip the cell to
# Text format just long enough to assign the literal string, then restore
# the cell to the default "Normal" style so no stray number formatting is
# left behind.
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "70.594.38"
$ws.Range("E2").Value = "  +2.17%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "3.539.91"
$ws.Range("E3").Value = "  +1.25%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.11%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "608.37"
$ws.Range("E5").Value = "  +4.81%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "173.54"
$ws.Range("E6").Value = "  +0.51%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.32%  "

# Row 8 - LidoStakedEther
Set-TextValue $ws.Range("D8") "3.534.10"
$ws.Range("E8").Value = "  +1.33%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  -0.12%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("D10") "0.200"
$ws.Range("E10").Value = "  +6.71%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +1.20%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -1.17%  "

# Row 13 - Avalanche
Set-TextValue $ws.Range("D13") "47.44"
$ws.Range("E13").Value = "  +1.88%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  +2.32%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "4.099.99"
$ws.Range("E15").Value = "  +0.98%  "

# Row 16 - BitcoinCash
Set-TextValue $ws.Range("D16") "627.11"
$ws.Range("E16").Value = "  -6.79%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  -2.73%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "70.593.96"
$ws.Range("E18").Value = "  +2.16%  "

# Row 19 - WrappedEther
Set-TextValue $ws.Range("D19") "3.536.83"
$ws.Range("E19").Value = "  +0.96%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  -1.70%  "

# Row 21 - Chainlink
Set-TextValue $ws.Range("D21") "17.43"
$ws.Range("E21").Value = "  +0.46%  "

# Row 22 - Polygon
Set-TextValue $ws.Range("D22") "0.887"
$ws.Range("E22").Value = "  -1.17%  "

# Row 23 - Uniswap
Set-TextValue $ws.Range("D23") "9.94"
$ws.Range("E23").Value = "  -10.58%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D24") "15.98"
$ws.Range("E24").Value = "  -0.32%  "

# Row 25 - Litecoin
Set-TextValue $ws.Range("D25") "97.06"
$ws.Range("E25").Value = "  -0.23%  "

# Row 26 - PancakeSwap
$ws.Range("E26").Value = "  +0.10%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  +0.14%  "

# Row 28 - ImmutableX
$ws.Range("E28").Value = "  -0.86%  "

# Row 29 - RenderToken
$ws.Range("E29").Value = "  -1.47%  "

# Row 30 - EthereumClassic
Set-TextValue $ws.Range("D30") "33.49"
$ws.Range("E30").Value = "  +2.13%  "

# Row 31 - Filecoin
Set-TextValue $ws.Range("D31") "8.48"
$ws.Range("E31").Value = "  -1.95%  "

# Row 32 - Stacks
$ws.Range("E32").Value = "  -2.27%  "

# Row 33 - Mantle
$ws.Range("E33").Value = "  -1.04%  "

# Row 34 - NEARProtocol
Set-TextValue $ws.Range("D34") "7.02"
$ws.Range("E34").Value = "  -2.53%  "

# Row 35 - Bittensor
Set-TextValue $ws.Range("D35") "568.11"
$ws.Range("E35").Value = "  -4.50%  "

# Row 36 - dogwifhat
Set-TextValue $ws.Range("D36") "3.65"
$ws.Range("E36").Value = "  +2.29%  "

# Row 37 - Cosmos
Set-TextValue $ws.Range("D37") "10.78"
$ws.Range("E37").Value = "  -0.38%  "

# Row 38 - OKB
Set-TextValue $ws.Range("D38") "57.50"
$ws.Range("E38").Value = "  +1.04%  "

# Row 39 - Hedera
$ws.Range("E39").Value = "  -1.49%  "

# Row 40 - FirstDigitalUSD
$ws.Range("E40").Value = "  +0.09%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  +6.04%  "

# Row 42 - VeChain
Set-TextValue $ws.Range("D42") "0.0456"
$ws.Range("E42").Value = "  +4.65%  "

# Row 43 - TheGraph
$ws.Range("E43").Value = "  -1.42%  "

# Row 44 - Maker
Set-TextValue $ws.Range("D44") "3.334.83"
$ws.Range("E44").Value = "  -2.09%  "

# Row 45 - ThetaToken
Set-TextValue $ws.Range("D45") "3.03"
$ws.Range("E45").Value = "  +5.39%  "

# Row 46 - PEPE
$ws.Range("E46").Value = "  +1.70%  "

# Row 47 - InjectiveProtocol
Set-TextValue $ws.Range("D47") "33.12"
$ws.Range("E47").Value = "  +0.02%  "

# Row 48 - Fetch.AI
$ws.Range("E48").Value = "  +2.90%  "

# Row 49 - Stellar
Set-TextValue $ws.Range("D49") "0.129"
$ws.Range("E49").Value = "  -2.04%  "

# Rows 50/51 - coin ranking order flips: Monero and MXToken trade places,
# with refreshed price/volume figures for both.
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D50") "5.74"
$ws.Range("E50").Value = "  +0.02%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D51") "133.58"
$ws.Range("E51").Value = "  +0.50%  "
